{"js": "// The document contains a single paragraph whose only content is a\n// bookmark (named \"_Hlk147415097\") wrapping one space character. The\n// commit \"upload of completed homework\" clears that placeholder\n// paragraph out, leaving a plain empty paragraph (no bookmark, no run,\n// no text) in its place.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Clear every paragraph's content (text + bookmarks). Deleting a\n// paragraph that is the sole/last paragraph of the body simply empties\n// it in place (Word always needs at least one paragraph), which is\n// exactly the \"<w:p/>\" result we want.\nfor (const paragraph of paragraphs.items) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single paragraph whose only content is a\n# bookmark (named \"_Hlk147415097\") wrapping one space character. The\n# commit \"upload of completed homework\" clears that placeholder\n# paragraph out, leaving a plain empty paragraph (no bookmark, no run,\n# no text) in its place.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    # Deleting the paragraph's range (not the paragraph mark itself)\n    # removes the run/text and any bookmarks anchored inside it, while\n    # leaving the paragraph mark (and thus the paragraph) in place.\n    $p.Range.Delete()\n}\n"}
